# Locate the target paragraph: the "Таблица 2." caption paragraph whose
# second run holds the sentence
#   "Тест-кейс для программы «Калькулятор КАСКО», функция расчета стоимости."
# (There is a near-identical sentence under "Таблица 1.", so we search only
# within the paragraph that also contains "Таблица 2.")

$d = $word.ActiveDocument

$oldSentence = "Тест-кейс для программы «Калькулятор КАСКО», функция расчета стоимости."

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like "*Таблица 2.*") -and ($t -like "*$oldSentence*")) {
        $target = $p
        break
    }
}

$full = $target.Range

# Find the exact sentence inside that paragraph and capture its Range.
$search = $full.Duplicate
$found = $search.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sentenceRange = $d.Range($search.Start, $search.End)

# Replace that single run with three runs that share the same run
# formatting (Times New Roman, sz 28 / szCs 28):
#   1) the sentence minus its trailing period (keeps rsidRPr="00F50845")
#   2) " в файл"
#   3) "."
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
          '<w:p>' +
          '<w:r w:rsidRPr="00F50845">' + $rPr + '<w:t>Тест-кейс для программы «Калькулятор КАСКО», функция расчета стоимости</w:t></w:r>' +
          '<w:r>' + $rPr + '<w:t xml:space="preserve"> в файл</w:t></w:r>' +
          '<w:r>' + $rPr + '<w:t>.</w:t></w:r>' +
          '</w:p>' +
          '</w:body>' +
          '</w:document>' +
          '</pkg:xmlData>' +
          '</pkg:part>' +
          '</pkg:package>'

$sentenceRange.InsertXML($newXml)
